# Generate Report for Handoff
# Updates the "latest generate/handoff" timestamps for the
# 90e65a5c-e474-4fac-ae46-5da719d0632e file row (row 7 on each sheet)
# to reflect a fresh handoff-report generation.

$wb = $excel.ActiveWorkbook

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-26 12:42:15"

# zh-cn sheet: column H = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-08-26 12:42:09"

# de-de sheet: column H = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-08-26 12:42:15"
